$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.611.99'
$ws.Range("E2").Value = '  -8.40%  '

$ws.Range("D3").Value = '1.653.73'
$ws.Range("E3").Value = '  -9.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.41%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.01'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -5.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5089'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -13.76%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.008'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2530'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -7.66%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.60'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -5.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06126'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -9.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07346'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.06%  '

$ws.Range("D12").Value = '1.648.14'
$ws.Range("E12").Value = '  -9.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.432'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5731'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -7.90%  '

$ws.Range("D15").Value = '1.877.49'
$ws.Range("E15").Value = '  -9.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008041'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -14.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.54'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -13.27%  '

$ws.Range("D18").Value = '26.615.00'
$ws.Range("E18").Value = '  -7.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.948'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -8.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.009'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.55'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -7.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '180.88'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -12.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.011'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.45%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.194'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -8.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.01'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -7.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.556'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1148'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -9.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '14.99'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -7.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.333'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.16%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05806'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -10.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.339'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -6.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.413'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -7.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.407'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -7.20%  '

$ws.Range("E34").Value = '  -6.17%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9747'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -7.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.430'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.20%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5926'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.637'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.87%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01573'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -7.42%  '

$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.069.77'
$ws.Range("E40").Value = '  -5.41%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8625'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.011'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.65%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.719'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -11.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '95.68'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.89%  '

$ws.Range("D45").Value = '1.793.29'
$ws.Range("E45").Value = '  -9.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000105'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -7.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.008'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.12'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -8.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4377'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.99%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05203'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -5.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.741'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.54%  '
